# Refresh the "maestro_precios" sheet with 5 newer days of FX data.
# The sheet is a rolling daily series (newest date on top, row 2) and
# keeps a fixed-size window of 310 data rows. We shift the existing
# data rows down by 5 (dropping the oldest 4 rows that fall out of the
# window) and write the 5 new top rows with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("maestro_precios")

$shiftBy = 5
$lastOldRow = 310        # old sheet's last data row (A1:C310)
$lastNewRow = 311        # new sheet's last data row (A1:C311)
$lastSourceRow = $lastNewRow - $shiftBy   # 306: last old row that still survives

# Shift surviving rows down, starting from the bottom so we never
# overwrite a source row before it has been read. Old rows below
# $lastSourceRow (i.e. 307..310) are oldest-dated and simply dropped.
for ($r = $lastSourceRow; $r -ge 2; $r--) {
    $destRow = $r + $shiftBy
    $aVal = $ws.Cells.Item($r, 1).Value2()
    $bVal = $ws.Cells.Item($r, 2).Value2()
    $cVal = $ws.Cells.Item($r, 3).Value2()

    $ws.Cells.Item($destRow, 1).Value = $aVal
    $ws.Cells.Item($destRow, 2).Value = $bVal
    $ws.Cells.Item($destRow, 3).Value = $cVal
}

# New top rows: maestro_id 22, newest dates/prices first.
$newTop = @(
    @(22, 46048, 1515.03),
    @(22, 46047, 1507.77),
    @(22, 46046, 1507.77),
    @(22, 46045, 1505.91),
    @(22, 46044, 1504.11)
)

$row = 2
foreach ($entry in $newTop) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Row $lastNewRow (311) previously fell outside the old used range
# (A1:C310), so it has no inherited number format yet. Match the date
# formatting used by the rest of column B.
$ws.Cells.Item($lastNewRow, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Clear whatever stale content still sits past the new last row so the
# sheet's used range/dimension shrinks back down correctly.
$extraFirst = $lastNewRow + 1
if ($extraFirst -le $lastOldRow) {
    $extraRange = $ws.Range("A" + $extraFirst + ":C" + $lastOldRow)
    $extraRange.ClearContents()
}

Write-Host "maestro_precios refreshed through row $lastNewRow"
